# Wed, Jul 15, 2020  1:05:08 AM
#
# 1) Table on slide 5 gets switched to a different (built-in) table style.
# 2) The deck's two theme parts swap identities: the theme that is wired to
#    the slide master (the one actually driving slide rendering) goes from
#    "Integral" / "Red Violet" colours to the stock "Office" colour
#    palette, while the notes-master theme moves the other way.
#    This host only exposes the slide-master-side theme for editing via
#    Slide.ThemeColorScheme, so we drive that one across to the Office
#    palette (font/format schemes in both theme parts are already
#    byte-identical, only the colour scheme differs).

$p = $ppt.ActivePresentation

# --- 1) Table style on slide 5 -------------------------------------------
$slide = $p.Slides.Item(5)
for ($i = 1; $i -le $slide.Shapes.Count; $i++) {
    $shp = $slide.Shapes.Item($i)
    if ($shp.HasTable) {
        $shp.Table.ApplyStyle("{0252FC5A-61D5-4A45-AA17-04E3762C93C8}")
    }
}

# --- 2) Theme colour scheme -> stock Office palette -----------------------
$tcs = $p.Slides.Item(1).ThemeColorScheme
$tcs.Colors(1).RGB  = 0        # dk1      000000
$tcs.Colors(2).RGB  = 16777215 # lt1      FFFFFF
$tcs.Colors(3).RGB  = 6968388  # dk2      44546A
$tcs.Colors(4).RGB  = 15132391 # lt2      E7E6E6
$tcs.Colors(5).RGB  = 13998939 # accent1  5B9BD5
$tcs.Colors(6).RGB  = 3243501  # accent2  ED7D31
$tcs.Colors(7).RGB  = 10855845 # accent3  A5A5A5
$tcs.Colors(8).RGB  = 49407    # accent4  FFC000
$tcs.Colors(9).RGB  = 12874308 # accent5  4472C4
$tcs.Colors(10).RGB = 4697456  # accent6  70AD47
$tcs.Colors(11).RGB = 12673797 # hlink    0563C1
$tcs.Colors(12).RGB = 7491477  # folHlink 954F72
